$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "X" status marker to the "Review" column (D) for rows 4, 5, 6
# Row 4: was in C4, now in D4
$ws.Range("C4").Value = $null
$ws.Range("D4").Value = "X"

# Row 5: was in B5, now in D5
$ws.Range("B5").Value = $null
$ws.Range("D5").Value = "X"

# Row 6: was in B6, now in D6
$ws.Range("B6").Value = $null
$ws.Range("D6").Value = "X"

# Update the active cell selection
$ws.Range("B6").Select()
